# Auto-generated edit script: updates ALC/ARM/BSM/CRP/CUL/LTW sheets
# to refresh currentAveragePrice / Leve profit data (scheduled runner sync).
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 2751.75
$ws.Range("I45").Value = 1007
$ws.Range("J45").Value = 3333.3333
$ws.Range("K45").Value = 3021
$ws.Range("L45").Value = 9999.999899999999
$ws.Range("M45").Value = -2829
$ws.Range("N45").Value = -10383.9999
$ws.Range("H137").Value = 13900719
$ws.Range("I137").Value = 893.6286
$ws.Range("J137").Value = 62550108
$ws.Range("K137").Value = 2680.8858
$ws.Range("L137").Value = 187650324
$ws.Range("M137").Value = -130.8858
$ws.Range("N137").Value = -187655424
$ws.Range("H138").Value = 2158.2537
$ws.Range("I138").Value = 1532.7142
$ws.Range("J138").Value = 3861.111
$ws.Range("K138").Value = 4598.142599999999
$ws.Range("L138").Value = 11583.333
$ws.Range("M138").Value = 541.8574000000008
$ws.Range("N138").Value = -21863.333

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1795221
$ws.Range("I61").Value = 801936.5600000001
$ws.Range("J61").Value = 14707920
$ws.Range("K61").Value = 801936.5600000001
$ws.Range("L61").Value = 14707920
$ws.Range("M61").Value = -801724.5600000001
$ws.Range("N61").Value = -14708344
$ws.Range("H74").Value = 32625356
$ws.Range("I74").Value = 25641666
$ws.Range("J74").Value = 66670844
$ws.Range("K74").Value = 25641666
$ws.Range("L74").Value = 66670844
$ws.Range("M74").Value = -25640792
$ws.Range("N74").Value = -66672592
$ws.Range("H77").Value = 32625356
$ws.Range("I77").Value = 25641666
$ws.Range("J77").Value = 66670844
$ws.Range("K77").Value = 128208330
$ws.Range("L77").Value = 333354220
$ws.Range("M77").Value = -128203962
$ws.Range("N77").Value = -333362956
$ws.Range("H88").Value = 4324.9443
$ws.Range("I88").Value = 2250
$ws.Range("J88").Value = 4917.7856
$ws.Range("K88").Value = 2250
$ws.Range("L88").Value = 4917.7856
$ws.Range("M88").Value = -1844
$ws.Range("N88").Value = -5729.7856
$ws.Range("H91").Value = 4324.9443
$ws.Range("I91").Value = 2250
$ws.Range("J91").Value = 4917.7856
$ws.Range("K91").Value = 2250
$ws.Range("L91").Value = 4917.7856
$ws.Range("M91").Value = -846
$ws.Range("N91").Value = -7725.7856
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 1795221
$ws.Range("I136").Value = 801936.5600000001
$ws.Range("J136").Value = 14707920
$ws.Range("K136").Value = 2405809.68
$ws.Range("L136").Value = 44123760
$ws.Range("M136").Value = -2403259.68
$ws.Range("N136").Value = -44128860

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 115986.664
$ws.Range("J87").Value = 115986.664
$ws.Range("L87").Value = 115986.664
$ws.Range("N87").Value = -118482.664
$ws.Range("H90").Value = 115986.664
$ws.Range("J90").Value = 115986.664
$ws.Range("L90").Value = 347959.992
$ws.Range("N90").Value = -360439.992
$ws.Range("H128").Value = 1783.3334
$ws.Range("I128").Value = 1783.3334
$ws.Range("K128").Value = 5350.0002
$ws.Range("M128").Value = -2860.0002

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1491430.9
$ws.Range("I31").Value = 1502.5151
$ws.Range("J31").Value = 6954501.5
$ws.Range("K31").Value = 1502.5151
$ws.Range("L31").Value = 6954501.5
$ws.Range("M31").Value = -1207.5151
$ws.Range("N31").Value = -6955091.5
$ws.Range("H34").Value = 1491430.9
$ws.Range("I34").Value = 1502.5151
$ws.Range("J34").Value = 6954501.5
$ws.Range("K34").Value = 1502.5151
$ws.Range("L34").Value = 6954501.5
$ws.Range("M34").Value = -1300.5151
$ws.Range("N34").Value = -6954905.5
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H107").Value = 455.6875
$ws.Range("I107").Value = 141.85
$ws.Range("K107").Value = 141.85
$ws.Range("M107").Value = 1778.15
$ws.Range("H134").Value = 1181396.2
$ws.Range("I134").Value = 4601.9653
$ws.Range("J134").Value = 8006803
$ws.Range("K134").Value = 13805.8959
$ws.Range("L134").Value = 24020409
$ws.Range("M134").Value = -11270.8959
$ws.Range("N134").Value = -24025479

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5251311.5
$ws.Range("I5").Value = 10989821
$ws.Range("J5").Value = 2382056.5
$ws.Range("K5").Value = 32969463
$ws.Range("L5").Value = 7146169.5
$ws.Range("M5").Value = -32969351
$ws.Range("N5").Value = -7146393.5
$ws.Range("H80").Value = 17383
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 17383
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 52149
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -54021
$ws.Range("H83").Value = 17383
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 17383
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 156447
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -165807
$ws.Range("H122").Value = 690.6667
$ws.Range("J122").Value = 1770.75
$ws.Range("L122").Value = 15936.75
$ws.Range("N122").Value = -20836.75
$ws.Range("H131").Value = 956.5909
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 1099.7059
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 3299.1177
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -13379.1177
$ws.Range("H135").Value = 5251311.5
$ws.Range("I135").Value = 10989821
$ws.Range("J135").Value = 2382056.5
$ws.Range("K135").Value = 98908389
$ws.Range("L135").Value = 21438508.5
$ws.Range("M135").Value = -98905854
$ws.Range("N135").Value = -21443578.5
$ws.Range("H140").Value = 4715.576
$ws.Range("I140").Value = 8305.117
$ws.Range("J140").Value = 3262.6667
$ws.Range("K140").Value = 24915.351
$ws.Range("L140").Value = 9788.000100000001
$ws.Range("M140").Value = -19735.351
$ws.Range("N140").Value = -20148.0001
$ws.Range("H141").Value = 3606
$ws.Range("I141").Value = 3507.5
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 10522.5
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -5342.5
$ws.Range("N141").Value = -22360

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 1400
$ws.Range("I107").Value = 1400
$ws.Range("K107").Value = 1400
$ws.Range("M107").Value = 520
